$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting the existing rows 6-7 down to 7-8
$ws.Rows.Item(6).Insert()

# Populate the newly-inserted row 6 with the new entry (Juveniles / caballeros / Vega, Ramón Emanuel)
$ws.Cells.Item(6, 1).Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Cells.Item(6, 2).Value = "Juveniles"
$ws.Cells.Item(6, 3).Value = "caballeros"
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = "Vega, Ramón Emanuel"
$ws.Cells.Item(6, 6).Value = 86
$ws.Cells.Item(6, 8).Value = 86

# G column stays blank for every data row (present-but-empty cell); copy that
# shape from a neighbouring row instead of assigning "" (which drops the cell).
$ws.Range("G5").Copy($ws.Cells.Item(6, 7))
